$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = "'29.491.90"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value2 = "  +1.97%  "
$ws.Cells.Item(3,4).Value2 = "'1.985.56"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value2 = "  +5.67%  "
$ws.Cells.Item(4,5).Value2 = "  -0.20%  "
$ws.Cells.Item(5,4).Value2 = "'326.02"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value2 = "  +0.33%  "
$ws.Cells.Item(6,4).Value2 = "'0.9997"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value2 = "  -0.23%  "
$ws.Cells.Item(7,4).Value2 = "'0.4687"
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value2 = "  +1.92%  "
$ws.Cells.Item(8,4).Value2 = "'0.3930"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value2 = "  +1.38%  "
$ws.Cells.Item(9,2).Value2 = "OKB"
$ws.Cells.Item(9,3).Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(9,4).Value2 = "'46.41"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value2 = "  -0.22%  "
$ws.Cells.Item(10,2).Value2 = "Dogecoin"
$ws.Cells.Item(10,3).Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10,4).Value2 = "'0.07936"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value2 = "  +0.83%  "
$ws.Cells.Item(11,2).Value2 = "Polygon"
$ws.Cells.Item(11,3).Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(11,4).Value2 = "'0.9996"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value2 = "  +1.37%  "
$ws.Cells.Item(12,2).Value2 = "Solana"
$ws.Cells.Item(12,3).Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(12,4).Value2 = "'22.90"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value2 = "  +5.22%  "
$ws.Cells.Item(13,2).Value2 = "WrappedEther"
$ws.Cells.Item(13,3).Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13,4).Value2 = "'1.951.75"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value2 = "  +3.75%  "
$ws.Cells.Item(14,2).Value2 = "Chainlink"
$ws.Cells.Item(14,3).Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(14,4).Value2 = "'7.245"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value2 = "  +3.67%  "
$ws.Cells.Item(15,2).Value2 = "Polkadot"
$ws.Cells.Item(15,3).Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(15,4).Value2 = "'5.861"
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value2 = "  +3.77%  "
$ws.Cells.Item(16,2).Value2 = "TRON"
$ws.Cells.Item(16,3).Value2 = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(16,4).Value2 = "'0.07129"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value2 = "  +2.44%  "
$ws.Cells.Item(17,2).Value2 = "Litecoin"
$ws.Cells.Item(17,3).Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(17,4).Value2 = "'88.58"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value2 = "  +0.66%  "
$ws.Cells.Item(18,2).Value2 = "BinanceUSD"
$ws.Cells.Item(18,3).Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(18,4).Value2 = "'1.002"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value2 = "  -0.04%  "
$ws.Cells.Item(19,2).Value2 = "ShibaInu"
$ws.Cells.Item(19,3).Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(19,4).Value2 = "'0.000009952"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value2 = "  -0.23%  "
$ws.Cells.Item(20,2).Value2 = "Avalanche"
$ws.Cells.Item(20,3).Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(20,4).Value2 = "'17.34"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value2 = "  +2.18%  "
$ws.Cells.Item(21,2).Value2 = "Dai"
$ws.Cells.Item(21,3).Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(21,4).Value2 = "'0.9990"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value2 = "  -0.20%  "
$ws.Cells.Item(22,2).Value2 = "WrappedBTC"
$ws.Cells.Item(22,3).Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(22,4).Value2 = "'29.590.98"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value2 = "  +2.28%  "
$ws.Cells.Item(23,2).Value2 = "Uniswap"
$ws.Cells.Item(23,3).Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(23,4).Value2 = "'5.526"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value2 = "  +5.40%  "
$ws.Cells.Item(24,2).Value2 = "Cosmos"
$ws.Cells.Item(24,3).Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(24,4).Value2 = "'11.26"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value2 = "  +2.74%  "
$ws.Cells.Item(25,4).Value2 = "'2.101"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value2 = "  +0.05%  "
$ws.Cells.Item(26,4).Value2 = "'157.72"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value2 = "  +0.93%  "
$ws.Cells.Item(27,4).Value2 = "'19.61"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value2 = "  +1.54%  "
$ws.Cells.Item(28,4).Value2 = "'5.970"
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value2 = "  -1.31%  "
$ws.Cells.Item(29,4).Value2 = "'120.11"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value2 = "  +2.32%  "
$ws.Cells.Item(30,4).Value2 = "'1.962"
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value2 = "  +1.77%  "
$ws.Cells.Item(31,4).Value2 = "'0.09435"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value2 = "  +1.20%  "
$ws.Cells.Item(32,4).Value2 = "'0.9047"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value2 = "  +0.26%  "
$ws.Cells.Item(33,4).Value2 = "'5.267"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value2 = "  +0.27%  "
$ws.Cells.Item(34,4).Value2 = "'1.347"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value2 = "  +2.22%  "
$ws.Cells.Item(35,4).Value2 = "'3.178"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value2 = "  -2.40%  "
$ws.Cells.Item(36,4).Value2 = "'0.05835"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value2 = "  +1.43%  "
$ws.Cells.Item(37,2).Value2 = "PEPE"
$ws.Cells.Item(37,3).Value2 = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(37,4).Value2 = "'0.000003416"
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value2 = "  +110.05%  "
$ws.Cells.Item(38,2).Value2 = "TrustWalletToken"
$ws.Cells.Item(38,3).Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38,4).Value2 = "'1.174"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value2 = "  -0.67%  "
$ws.Cells.Item(39,4).Value2 = "'0.02117"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value2 = "  +2.29%  "
$ws.Cells.Item(40,4).Value2 = "'7.899"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value2 = "  +2.97%  "
$ws.Cells.Item(41,4).Value2 = "'0.5753"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value2 = "  +1.80%  "
$ws.Cells.Item(42,5).Value2 = "  +3.45%  "
$ws.Cells.Item(43,4).Value2 = "'9.803"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value2 = "  +1.34%  "
$ws.Cells.Item(44,4).Value2 = "'11.99"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value2 = "  +0.91%  "
$ws.Cells.Item(45,4).Value2 = "'0.5367"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value2 = "  +0.30%  "
$ws.Cells.Item(46,5).Value2 = "  +6.19%  "
$ws.Cells.Item(47,4).Value2 = "'2.176"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value2 = "  -4.11%  "
$ws.Cells.Item(49,5).Value2 = "  +1.17%  "
$ws.Cells.Item(50,4).Value2 = "'114.30"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value2 = "  +1.34%  "
$ws.Cells.Item(51,4).Value2 = "'0.3099"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value2 = "  +8.24%  "
